# EPBDS-11327 "class" information can't be found in the OpenL datatypes
#
# Adds new test cases (case21..case24) to the SRClass field-access test
# sheet, covering `.class.simpleName` access on a bean instance and on a
# user-defined datatype (MyType), plus a couple of supporting rows that
# exercise direct field / getter access on the new MyType datatype.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 / 27: fill in the secondary "mini" test table (cols E:H) for
#     the two existing steps (sr3:ChildWithStaticId / case20) that were
#     previously left blank, plus a brand-new step (case21) testing
#     AccessBean.class.simpleName.
$ws.Range("E26").Value = '_res_.$Value$case21'
$ws.Range("F26").Value = 'SR21'
$ws.Range("G26").Value = 'AccessBean'
$ws.Range("H26").Value = 'AccessBean'

$ws.Range("E27").Value = '_res_.$Value$case22'
$ws.Range("F27").Value = 'SR22'
$ws.Range("G27").Value = 'MyType'
$ws.Range("H27").Value = 'MyType'

# --- Row 28: new step case21 - AccessBean.class.simpleName
$ws.Range("B28").Value = 'case21'
$ws.Range("C28").Value = "'= AccessBean.class.simpleName"
$ws.Range("E28").Value = '_res_.$Value$case23'
$ws.Range("F28").Value = 'SR23'
$ws.Range("G28").Value = 'XXX'
$ws.Range("H28").Value = 'XXX'

# --- Row 29: new step case22 - MyType.class.simpleName
$ws.Range("B29").Value = 'case22'
$ws.Range("C29").Value = "'= MyType.class.simpleName"
$ws.Range("E29").Value = '_res_.$Value$case24'
$ws.Range("F29").Value = 'SR24'
$ws.Range("G29").Value = 'XXX'
$ws.Range("H29").Value = 'XXX'

# --- Row 30: new step case23 - new MyType().value (direct field access)
$ws.Range("B30").Value = 'case23'
$ws.Range("C30").Value = "'= new MyType().value"

# --- Row 31: new step case24 - new MyType().getValue() (getter access)
$ws.Range("B31").Value = 'case24'
$ws.Range("C31").Value = "'= new MyType().getValue()"

# Apply the same thin all-around border already used by the rest of the
# E:H mini-table to the newly populated range.
$ws.Range("E26:H29").Borders.LineStyle = 1
$ws.Range("E26:H29").Borders.Weight = 2

# --- Rows 35/36: new "Datatype MyType" block documenting the datatype
#     used by the new test cases (MyType { String value = "XXX" }).
$ws.Range("B35").Value = 'Datatype MyType'
$ws.Range("B36").Value = 'String'
$ws.Range("C36").Value = 'value'
$ws.Range("D36").Value = 'XXX'

Write-Output "EPBDS-11327 edit applied"
